$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Revert "attempt to add TOC filtering" by removing the TOC Filter / All TOCs
# row that was inserted as row 2, shifting the remaining data back up.
$ws.Rows.Item(2).Delete()

$ws.Range("B2").Select()
